$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "done"
$ws.Range("E3").Value = "done"

$ws.Range("B2").Select()
